# Apply updated loading_percent values (case with 380 kV) to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 11.06403580291998
$ws.Range("C2").Value = 4.377007299558759
$ws.Range("D2").Value = 14.94406876482549
$ws.Range("E2").Value = 16.35678265942679
$ws.Range("G2").Value = 34.77004468297362
$ws.Range("H2").Value = 15.83897694954581
$ws.Range("J2").Value = 9.334096853267642
$ws.Range("K2").Value = 10.42670214008421
$ws.Range("N2").Value = 19.12543094865401
$ws.Range("O2").Value = 24.90485738767283
# Row 3
$ws.Range("B3").Value = 10.80301688503346
$ws.Range("C3").Value = 4.165957524377777
$ws.Range("D3").Value = 14.87956042906489
$ws.Range("E3").Value = 16.29229328350535
$ws.Range("G3").Value = 34.82110890560872
$ws.Range("H3").Value = 15.88369070356707
$ws.Range("J3").Value = 9.340605682173548
$ws.Range("K3").Value = 10.24874268582078
$ws.Range("N3").Value = 19.18491144972413
$ws.Range("O3").Value = 24.9719954470897
# Row 4
$ws.Range("B4").Value = 10.64155852741815
$ws.Range("C4").Value = 4.029951350874351
$ws.Range("D4").Value = 14.84316532035108
$ws.Range("E4").Value = 16.25624913527222
$ws.Range("G4").Value = 34.8620878796265
$ws.Range("H4").Value = 15.91353829265657
$ws.Range("J4").Value = 9.34598782413817
$ws.Range("K4").Value = 10.1395938675004
$ws.Range("N4").Value = 19.22313279394131
$ws.Range("O4").Value = 25.01811877923332
# Row 5
$ws.Range("B5").Value = 10.57556235041867
$ws.Range("C5").Value = 3.972953084123236
$ws.Range("D5").Value = 14.82915385038628
$ws.Range("E5").Value = 16.24246555571386
$ws.Range("G5").Value = 34.88120109367515
$ws.Range("H5").Value = 15.92630311064896
$ws.Range("J5").Value = 9.348529875553
$ws.Range("K5").Value = 10.0952014854398
$ws.Range("N5").Value = 19.23913708203697
$ws.Range("O5").Value = 25.03814459938284
# Row 6
$ws.Range("B6").Value = 10.56459458946149
$ws.Range("C6").Value = 3.963394768307788
$ws.Range("D6").Value = 14.82687709678217
$ws.Range("E6").Value = 16.24023175479338
$ws.Range("G6").Value = 34.88452041634758
$ws.Range("H6").Value = 15.9284590361634
$ws.Range("J6").Value = 9.348973055330207
$ws.Range("K6").Value = 10.08783706733119
$ws.Range("N6").Value = 19.24182051686411
$ws.Range("O6").Value = 25.04154411016138
# Row 7
$ws.Range("B7").Value = 10.64066915837942
$ws.Range("C7").Value = 4.029188970031348
$ws.Range("D7").Value = 14.84297302247006
$ws.Range("E7").Value = 16.25605956779932
$ws.Range("G7").Value = 34.86233588243325
$ws.Range("H7").Value = 15.91370800740704
$ws.Range("J7").Value = 9.346020694526592
$ws.Range("K7").Value = 10.13899475007745
$ws.Range("N7").Value = 19.2233468955581
$ws.Range("O7").Value = 25.01838387625402
# Row 8
$ws.Range("B8").Value = 10.97434203588207
$ws.Range("C8").Value = 4.305591766609813
$ws.Range("D8").Value = 14.92116679596398
$ws.Range("E8").Value = 16.33381654677598
$ws.Range("G8").Value = 34.78565034458819
$ws.Range("H8").Value = 15.85389749574218
$ws.Range("J8").Value = 9.336053723530444
$ws.Range("K8").Value = 10.36535082480615
$ws.Range("N8").Value = 19.14558774986147
$ws.Range("O8").Value = 24.92698834927973
# Row 9
$ws.Range("B9").Value = 11.61507664271142
$ws.Range("C9").Value = 4.795251911110246
$ws.Range("D9").Value = 15.0994320857208
$ws.Range("E9").Value = 16.51393521307705
$ws.Range("G9").Value = 34.71188418810839
$ws.Range("H9").Value = 15.75560217032308
$ws.Range("J9").Value = 9.327486036600106
$ws.Range("K9").Value = 10.80773116230705
$ws.Range("N9").Value = 19.00653214049431
$ws.Range("O9").Value = 24.7867305129738
# Row 10
$ws.Range("B10").Value = 12.07232836578513
$ws.Range("C10").Value = 5.121576712459789
$ws.Range("D10").Value = 15.24477214509719
$ws.Range("E10").Value = 16.66230872535537
$ws.Range("G10").Value = 34.70465536607496
$ws.Range("H10").Value = 15.69496873900329
$ws.Range("J10").Value = 9.327856648749879
$ws.Range("K10").Value = 11.12864537625443
$ws.Range("N10").Value = 18.91247046009035
$ws.Range("O10").Value = 24.70755219758524
# Row 11
$ws.Range("B11").Value = 12.27641844868482
$ws.Range("C11").Value = 5.262546494258237
$ws.Range("D11").Value = 15.31381551387775
$ws.Range("E11").Value = 16.73309776789767
$ws.Range("G11").Value = 34.71159617492079
$ws.Range("H11").Value = 15.66990149270735
$ws.Range("J11").Value = 9.329465047433153
$ws.Range("K11").Value = 11.2731063604596
$ws.Range("N11").Value = 18.87142065965516
$ws.Range("O11").Value = 24.67673630501803
# Row 12
$ws.Range("B12").Value = 12.35306372380256
$ws.Range("C12").Value = 5.314838788305674
$ws.Range("D12").Value = 15.34036366164677
$ws.Range("E12").Value = 16.76035963343499
$ws.Range("G12").Value = 34.71569576037309
$ws.Range("H12").Value = 15.6607709661093
$ws.Range("J12").Value = 9.330280330736818
$ws.Range("K12").Value = 11.32754023673302
$ws.Range("N12").Value = 18.8561249470402
$ws.Range("O12").Value = 24.6658167755737
# Row 13
$ws.Range("B13").Value = 12.33658637109288
$ws.Range("C13").Value = 5.303625383285273
$ws.Range("D13").Value = 15.3346284177102
$ws.Range("E13").Value = 16.75446833484894
$ws.Range("G13").Value = 34.7147474163392
$ws.Range("H13").Value = 15.66272128974597
$ws.Range("J13").Value = 9.330095587475629
$ws.Range("K13").Value = 11.31582974005172
$ws.Range("N13").Value = 18.85940809836266
$ws.Range("O13").Value = 24.66813512071383
# Row 14
$ws.Range("B14").Value = 12.28273732744629
$ws.Range("C14").Value = 5.266870529611306
$ws.Range("D14").Value = 15.31599167666483
$ws.Range("E14").Value = 16.73533159298083
$ws.Range("G14").Value = 34.7119039672123
$ws.Range("H14").Value = 15.66914306469826
$ws.Range("J14").Value = 9.329527993442809
$ws.Range("K14").Value = 11.27759037087852
$ws.Range("N14").Value = 18.8701572893181
$ws.Range("O14").Value = 24.67582291328629
# Row 15
$ws.Range("B15").Value = 12.24966783053676
$ws.Range("C15").Value = 5.244214797967168
$ws.Range("D15").Value = 15.30462805943833
$ws.Range("E15").Value = 16.72366857065632
$ws.Range("G15").Value = 34.7103538591639
$ws.Range("H15").Value = 15.67312372282166
$ws.Range("J15").Value = 9.329207155192377
$ws.Range("K15").Value = 11.25413096797453
$ws.Range("N15").Value = 18.87677386639115
$ws.Range("O15").Value = 24.68062959736421
# Row 16
$ws.Range("B16").Value = 12.0589051610443
$ws.Range("C16").Value = 5.112212220291627
$ws.Range("D16").Value = 15.24031738069811
$ws.Range("E16").Value = 16.65774724515859
$ws.Range("G16").Value = 34.7044076621998
$ws.Range("H16").Value = 15.69665757330378
$ws.Range("J16").Value = 9.327780440414632
$ws.Range("K16").Value = 11.11916940699793
$ws.Range("N16").Value = 18.9151880616485
$ws.Range("O16").Value = 24.70967092434662
# Row 17
$ws.Range("B17").Value = 11.94081946729231
$ws.Range("C17").Value = 5.029306401471159
$ws.Range("D17").Value = 15.20160204959573
$ws.Range("E17").Value = 16.61813759057416
$ws.Range("G17").Value = 34.7033804494757
$ws.Range("H17").Value = 15.71173918796955
$ws.Range("J17").Value = 9.327273410036957
$ws.Range("K17").Value = 11.03594724528913
$ws.Range("N17").Value = 18.93919858320938
$ws.Range("O17").Value = 24.72882059655693
# Row 18
$ws.Range("B18").Value = 11.87253528542708
$ws.Range("C18").Value = 4.980918060444007
$ws.Range("D18").Value = 15.17961068896438
$ws.Range("E18").Value = 16.59566626480689
$ws.Range("G18").Value = 34.70375251353882
$ws.Range("H18").Value = 15.72065046500503
$ws.Range("J18").Value = 9.327117356377551
$ws.Range("K18").Value = 10.98793947605606
$ws.Range("N18").Value = 18.95317258258889
$ws.Range("O18").Value = 24.74032459884853
# Row 19
$ws.Range("B19").Value = 11.84935526024548
$ws.Range("C19").Value = 4.964414339949148
$ws.Range("D19").Value = 15.17221284881573
$ws.Range("E19").Value = 16.58811182570409
$ws.Range("G19").Value = 34.70404383446846
$ws.Range("H19").Value = 15.72370832243405
$ws.Range("J19").Value = 9.327087828265505
$ws.Range("K19").Value = 10.97166232304508
$ws.Range("N19").Value = 18.9579321039255
$ws.Range("O19").Value = 24.74430369602869
# Row 20
$ws.Range("B20").Value = 11.95342817613959
$ws.Range("C20").Value = 5.038204752668977
$ws.Range("D20").Value = 15.20569485637086
$ws.Range("E20").Value = 16.62232203582136
$ws.Range("G20").Value = 34.70339013614628
$ws.Range("H20").Value = 15.71010922196429
$ws.Range("J20").Value = 9.327313358260685
$ws.Range("K20").Value = 11.04482131448446
$ws.Range("N20").Value = 18.93662567974217
$ws.Range("O20").Value = 24.72673139529932
# Row 21
$ws.Range("B21").Value = 12.29857202237821
$ws.Range("C21").Value = 5.277696003061118
$ws.Range("D21").Value = 15.32145495136717
$ws.Range("E21").Value = 16.74094030044575
$ws.Range("G21").Value = 34.71269923199246
$ws.Range("H21").Value = 15.66724701036541
$ws.Range("J21").Value = 9.329689119937695
$ws.Range("K21").Value = 11.28882993347755
$ws.Range("N21").Value = 18.86699324394981
$ws.Range("O21").Value = 24.6735444612451
# Row 22
$ws.Range("B22").Value = 12.52038647953855
$ws.Range("C22").Value = 5.427859801650977
$ws.Range("D22").Value = 15.39945055510042
$ws.Range("E22").Value = 16.82111085240387
$ws.Range("G22").Value = 34.72735804313688
$ws.Range("H22").Value = 15.64134356994678
$ws.Range("J22").Value = 9.332443360168577
$ws.Range("K22").Value = 11.44670583369921
$ws.Range("N22").Value = 18.82293499791189
$ws.Range("O22").Value = 24.64315465507504
# Row 23
$ws.Range("B23").Value = 12.40236732959244
$ws.Range("C23").Value = 5.348300377386846
$ws.Range("D23").Value = 15.35761489096577
$ws.Range("E23").Value = 16.77808628705516
$ws.Range("G23").Value = 34.71875001196489
$ws.Range("H23").Value = 15.65497563844795
$ws.Range("J23").Value = 9.330863726756119
$ws.Range("K23").Value = 11.36260694920628
$ws.Range("N23").Value = 18.84631736487747
$ws.Range("O23").Value = 24.65897383040334
# Row 24
$ws.Range("B24").Value = 11.94772900379605
$ws.Range("C24").Value = 5.034184064194946
$ws.Range("D24").Value = 15.20384366664572
$ws.Range("E24").Value = 16.62042930984697
$ws.Range("G24").Value = 34.70338275831224
$ws.Range("H24").Value = 15.71084537996833
$ws.Range("J24").Value = 9.327294875711576
$ws.Range("K24").Value = 11.04080984921346
$ws.Range("N24").Value = 18.93778835966041
$ws.Range("O24").Value = 24.72767438264243
# Row 25
$ws.Range("B25").Value = 11.44375452052051
$ws.Range("C25").Value = 4.668566440661005
$ws.Range("D25").Value = 15.04861961052124
$ws.Range("E25").Value = 16.46232891133619
$ws.Range("G25").Value = 34.72360609055888
$ws.Range("H25").Value = 15.78015972745436
$ws.Range("J25").Value = 9.328631016598948
$ws.Range("K25").Value = 10.68855180356959
$ws.Range("N25").Value = 19.04272140734774
$ws.Range("O25").Value = 24.82048990599929
